# Rename the two logo pictures that live in the document's headers/footers.
#
#   * BTec_Logo-Orange  (word/header1.xml, word/header2.xml) : image1.jpg -> image2.jpg
#   * PearsonLogo.png   (word/footer1.xml, word/footer2.xml) : image2.png -> image1.png
#
# (The alt-text / description on each picture is unchanged - only the
# picture's Name changes.)

$d = $word.ActiveDocument
$sec = $d.Sections(1)

function Rename-InlineLogo($headerFooter, $newName) {
    $rng = $headerFooter.Range
    $count = $rng.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $para = $rng.Paragraphs($i)
        if ($para.Range.InlineShapes.Count -gt 0) {
            $para.Range.InlineShapes(1).Name = $newName
        }
    }
}

# Pearson logo in both footers: image2.png -> image1.png
Rename-InlineLogo $sec.Footers(1) "image1.png"
Rename-InlineLogo $sec.Footers(2) "image1.png"

# BTec logo in both headers: image1.jpg -> image2.jpg
Rename-InlineLogo $sec.Headers(1) "image2.jpg"
Rename-InlineLogo $sec.Headers(2) "image2.jpg"
